# SampleTestData.xlsx — "added tests for snapdeal"
#
# Sheet1 holds the Selenium search-box test data (S.NO / SearchData / Name).
# Swap the old sample search terms ("cooking", "obama") for the new
# snapdeal search terms ("mobile", "pen"); the Name column (Vivek / Arvind)
# stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "mobile"
$ws.Range("B3").Value = "pen"

# Leave the selection where it ended up after editing the data.
$null = $ws.Range("E14").Select()

$null = $wb.Save()
